$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inputs")

# Update the NOTE text in Y1: clarify how to find climate station IDs.
$ws.Range("Y1").Value = "NOTE:  Please use the  ""Lookup Tables"" tab to get the possible values for Units, Soil Texture, and Slope Shape.  To find the climate station IDs, visit the Climate Station section of the RHEM Web Tool: https://apps.tucson.ars.ag.gov/rhem"

# Update the Avg Precipitation (mm/year) value in S2, keeping it as literal
# text (leading spaces preserved) rather than letting it be coerced to a
# number. Format as Text first, set the value, then drop back to the
# workbook's default "Normal" style so no stray formatting is left behind.
$ws.Range("S2").NumberFormat = "@"
$ws.Range("S2").Value = "   261.94"
$ws.Range("S2").Style = "Normal"
